$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 (pushes existing rows 23-41 down to 24-42)
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly price record
$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value = "La Araucanía"
$ws.Cells.Item(23, 4).Value = 44729
$ws.Cells.Item(23, 5).Value = 9
$ws.Cells.Item(23, 6).Value = 100112010
$ws.Cells.Item(23, 7).Value = "Achicoria"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 65
$ws.Cells.Item(23, 11).Value = 8000
$ws.Cells.Item(23, 12).Value = 8000
$ws.Cells.Item(23, 13).Value = 8000
$ws.Cells.Item(23, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 444
$ws.Cells.Item(23, 17).Value = 18
$ws.Cells.Item(23, 18).Value = "Hortaliza"
